$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.236.08"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "'1.860.72"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("D4").Value = "'0.9993"

$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "'0.7049"
$ws.Range("E5").Value = "  -1.18%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'242.38"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.9994"

$ws.Range("D8").Value = "'0.3125"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "'0.07802"
$ws.Range("E9").Value = "  -3.03%  "

$ws.Range("D10").Value = "'24.25"
$ws.Range("E10").Value = "  -4.10%  "

$ws.Range("D11").Value = "'0.07998"
$ws.Range("E11").Value = "  -4.13%  "

$ws.Range("D12").Value = "'1.867.33"
$ws.Range("E12").Value = "  -1.79%  "

$ws.Range("D13").Value = "'94.18"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").Value = "'5.178"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").Value = "'0.6972"
$ws.Range("E15").Value = "  -3.08%  "

$ws.Range("D16").Value = "'6.395"
$ws.Range("E16").Value = "  +0.93%  "

$ws.Range("D17").Value = "'29.222.91"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("D18").Value = "'0.000008298"
$ws.Range("E18").Value = "  -3.37%  "

$ws.Range("D19").Value = "'253.67"
$ws.Range("E19").Value = "  +4.56%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'13.12"
$ws.Range("E20").Value = "  -1.05%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "'2.107.00"
$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").Value = "'7.526"
$ws.Range("E23").Value = "  -4.36%  "

$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").Value = "'0.1559"
$ws.Range("E25").Value = "  -2.15%  "

$ws.Range("D26").Value = "'8.988"
$ws.Range("E26").Value = "  -1.11%  "

$ws.Range("D27").Value = "'159.62"
$ws.Range("E27").Value = "  -2.43%  "

$ws.Range("D28").Value = "'18.91"
$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("D29").Value = "'1.499"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").Value = "'4.309"
$ws.Range("E30").Value = "  -2.54%  "

$ws.Range("D31").Value = "'4.266"
$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("D32").Value = "'1.210"
$ws.Range("E32").Value = "  +0.89%  "

$ws.Range("D33").Value = "'0.05274"
$ws.Range("E33").Value = "  -2.06%  "

$ws.Range("D34").Value = "'1.886"
$ws.Range("E34").Value = "  -3.27%  "

$ws.Range("D35").Value = "'0.7479"
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").Value = "'1.155"
$ws.Range("E36").Value = "  -2.36%  "

$ws.Range("D37").Value = "'2.711"
$ws.Range("E37").Value = "  +0.50%  "

$ws.Range("D38").Value = "'0.01868"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("D39").Value = "'1.245.97"
$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("D40").Value = "'2.733"
$ws.Range("E40").Value = "  -0.50%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8976"
$ws.Range("E41").Value = "  -2.16%  "

$ws.Range("D42").Value = "'110.97"
$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.132"
$ws.Range("E43").Value = "  -7.09%  "

$ws.Range("D44").Value = "'70.80"
$ws.Range("E44").Value = "  -5.65%  "

$ws.Range("D45").Value = "'0.9989"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").Value = "'2.007.37"
$ws.Range("E47").Value = "  -1.53%  "

$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("D49").Value = "'1.784"
$ws.Range("E49").Value = "  -1.61%  "

$ws.Range("D50").Value = "'9.469"
$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("D51").Value = "'0.4296"
$ws.Range("E51").Value = "  -2.24%  "
